# Update "想去人数" (want-to-go count) figures that changed between the
# previous data pull and the regenerated gh-pages output.
#
# Sheet "展览" (Exhibitions) and sheet "全部类型" (All types) share the
# same exhibition rows, and sheet "演出" (Shows) / "全部类型" share the
# same show row, so each numeric update is applied once per sheet that
# contains the corresponding row.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibitions) sheet updates
$wsExpo.Range("F4").Value  = 2042
$wsExpo.Range("F5").Value  = 336
$wsExpo.Range("F6").Value  = 582
$wsExpo.Range("F9").Value  = 10523
$wsExpo.Range("F11").Value = 153
$wsExpo.Range("F12").Value = 277
$wsExpo.Range("F14").Value = 408
$wsExpo.Range("F15").Value = 7419
$wsExpo.Range("F17").Value = 704
$wsExpo.Range("F18").Value = 203
$wsExpo.Range("F20").Value = 3305

# 演出 (Shows) sheet update
$wsShow.Range("F2").Value = 21

# 全部类型 (All types) sheet updates
$wsAll.Range("F4").Value  = 2042
$wsAll.Range("F5").Value  = 336
$wsAll.Range("F6").Value  = 582
$wsAll.Range("F7").Value  = 21
$wsAll.Range("F12").Value = 10523
$wsAll.Range("F14").Value = 153
$wsAll.Range("F15").Value = 277
$wsAll.Range("F17").Value = 408
$wsAll.Range("F18").Value = 7419
$wsAll.Range("F20").Value = 704
$wsAll.Range("F21").Value = 203
$wsAll.Range("F23").Value = 3305
